$d = $word.ActiveDocument

$replacements = @(
    @{old="644×4=2576"; new="862×9=7758"},
    @{old="383×7=2681"; new="646×9=5814"},
    @{old="775×8=6200"; new="756×6=4536"},
    @{old="949×5=4745"; new="144×3=432"},
    @{old="952×7=6664"; new="610×9=5490"},
    @{old="290×2=580"; new="725×8=5800"},
    @{old="789×2=1578"; new="908×3=2724"},
    @{old="284×2=568"; new="572×9=5148"},
    @{old="467×9=4203"; new="236×6=1416"},
    @{old="802×2=1604"; new="685×6=4110"},
    @{old="205×3=615"; new="978×3=2934"},
    @{old="813×4=3252"; new="138×2=276"},
    @{old="426×2=852"; new="197×4=788"},
    @{old="763×8=6104"; new="596×8=4768"},
    @{old="610×5=3050"; new="804×6=4824"},
    @{old="784×7=5488"; new="448×6=2688"},
    @{old="729×7=5103"; new="641×6=3846"},
    @{old="351×9=3159"; new="534×6=3204"},
    @{old="784×5=3920"; new="144×6=864"},
    @{old="291×3=873"; new="519×4=2076"},
    @{old="911×3=2733"; new="879×6=5274"},
    @{old="474×6=2844"; new="475×2=950"},
    @{old="551×7=3857"; new="258×2=516"},
    @{old="556×9=5004"; new="389×7=2723"},
    @{old="835×2=1670"; new="855×4=3420"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

Write-Output "Done applying $($replacements.Count) replacements"
